$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.175.67"
$ws.Range("E2").Value = "  -0.63%  "

$ws.Range("D3").Value = "1.860.97"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "0.7144"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D6").Value = "240.22"
$ws.Range("E6").Value = "  +0.66%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "0.07714"
$ws.Range("E8").Value = "  -1.44%  "

$ws.Range("D9").Value = "0.3071"
$ws.Range("E9").Value = "  +0.16%  "

$ws.Range("D10").Value = "24.93"
$ws.Range("E10").Value = "  -1.56%  "

$ws.Range("D11").Value = "0.08239"
$ws.Range("E11").Value = "  +0.66%  "

$ws.Range("D12").Value = "1.863.82"
$ws.Range("E12").Value = "  -0.89%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.7160"
$ws.Range("E13").Value = "  -0.87%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.212"
$ws.Range("E14").Value = "  -0.66%  "

$ws.Range("D15").Value = "90.32"
$ws.Range("E15").Value = "  +1.14%  "

$ws.Range("D16").Value = "29.161.53"
$ws.Range("E16").Value = "  -0.68%  "

$ws.Range("D17").Value = "5.849"
$ws.Range("E17").Value = "  +0.55%  "

$ws.Range("D18").Value = "243.34"
$ws.Range("E18").Value = "  +0.36%  "

$ws.Range("D19").Value = "0.000007784"
$ws.Range("E19").Value = "  -0.71%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "13.12"
$ws.Range("E20").Value = "  -1.46%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.109.29"
$ws.Range("E21").Value = "  -1.06%  "

$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").Value = "7.935"
$ws.Range("E23").Value = "  +2.19%  "

$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").Value = "0.1587"
$ws.Range("E25").Value = "  +7.88%  "

$ws.Range("D26").Value = "162.24"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").Value = "8.898"
$ws.Range("E27").Value = "  -0.75%  "

$ws.Range("D28").Value = "18.17"
$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "1.318"
$ws.Range("E29").Value = "  -3.28%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.492"
$ws.Range("E30").Value = "  +0.67%  "

$ws.Range("D31").Value = "4.344"
$ws.Range("E31").Value = "  +0.82%  "

$ws.Range("D32").Value = "4.082"
$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("D33").Value = "0.05188"
$ws.Range("E33").Value = "  -0.70%  "

$ws.Range("D34").Value = "1.909"
$ws.Range("E34").Value = "  -1.13%  "

$ws.Range("D35").Value = "1.173"
$ws.Range("E35").Value = "  -1.53%  "

$ws.Range("D36").Value = "0.7280"
$ws.Range("E36").Value = "  +1.14%  "

$ws.Range("D37").Value = "2.679"
$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").Value = "0.01846"
$ws.Range("E38").Value = "  -0.52%  "

$ws.Range("E39").Value = "  -0.13%  "

$ws.Range("D40").Value = "1.149.79"
$ws.Range("E40").Value = "  -1.91%  "

$ws.Range("D41").Value = "0.9009"
$ws.Range("E41").Value = "  -1.43%  "

$ws.Range("D42").Value = "6.079"
$ws.Range("E42").Value = "  +1.53%  "

$ws.Range("D43").Value = "72.29"
$ws.Range("E43").Value = "  +1.07%  "

$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("D45").Value = "101.51"
$ws.Range("E45").Value = "  -0.78%  "

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "2.007.77"
$ws.Range("E46").Value = "  -0.74%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.5266"
$ws.Range("E47").Value = "  -1.11%  "

$ws.Range("D48").Value = "1.763"
$ws.Range("E48").Value = "  -0.16%  "

$ws.Range("D49").Value = "9.270"
$ws.Range("E49").Value = "  +0.50%  "

$ws.Range("D50").Value = "2.870"
$ws.Range("E50").Value = "  -1.87%  "

$ws.Range("D51").Value = "0.9963"
$ws.Range("E51").Value = "  -0.82%  "
